# Update scripts with new TPM - recalculated receptor/edge specificity values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Target cluster = ECs)
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1204806666666667
$ws.Range("N2").Value = 0.361442
$ws.Range("O2").Value = 0.0341140812585293
$ws.Range("P2").Value = 0.03411408125852929
$ws.Range("Q2").Value = 0.05468320274355556
$ws.Range("R2").Value = 0.492148824692
$ws.Range("S2").Value = 0.0341140812585293
$ws.Range("T2").Value = 0.03411408125852929

# Row 3 (Target cluster = FAPs)
$ws.Range("O3").Value = 0.8134571113412176
$ws.Range("P3").Value = 0.8134571113412176
$ws.Range("S3").Value = 0.8134571113412176
$ws.Range("T3").Value = 0.8134571113412176

# Row 4 (Target cluster = MuSCs)
$ws.Range("M4").Value = 0.5383326666666667
$ws.Range("N4").Value = 1.614998
$ws.Range("O4").Value = 0.1524288074002532
$ws.Range("P4").Value = 0.1524288074002531
$ws.Range("Q4").Value = 0.2443359185275556
$ws.Range("R4").Value = 2.199023266748
$ws.Range("S4").Value = 0.1524288074002532
$ws.Range("T4").Value = 0.1524288074002531
